# Updated cryptos list: refresh Price (col D) and Volume(1h) (col E) text values.
# All source values are plain text cells (t="inlineStr" in the original OOXML),
# so assignments that Excel would otherwise auto-coerce into numbers (single-dot
# decimals like "596.81") are forced back to Text via NumberFormat="@" and then
# the style is reset to "Normal" so no stray style index is left on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.408.11'
$ws.Range("E2").Value = '  +0.75%  '
$ws.Range("D3").Value = '3.492.72'
$ws.Range("E3").Value = '  -0.20%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '596.81'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.39%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '179.36'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.89%  '
$ws.Range("E8").Value = '  +2.17%  '
$ws.Range("D9").Value = '3.496.70'
$ws.Range("E9").Value = '  -0.04%  '
$ws.Range("E10").Value = '  +4.80%  '
$ws.Range("E11").Value = '  -1.98%  '
$ws.Range("E12").Value = '  +1.12%  '
$ws.Range("D13").Value = '4.097.97'
$ws.Range("E13").Value = '  -0.11%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '32.39'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +10.61%  '
$ws.Range("E15").Value = '  +0.64%  '
$ws.Range("D16").Value = '67.402.14'
$ws.Range("E16").Value = '  +0.72%  '
$ws.Range("E17").Value = '  -0.66%  '
$ws.Range("D18").Value = '3.487.98'
$ws.Range("E18").Value = '  -0.29%  '
$ws.Range("E19").Value = '  -0.02%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.31'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.19%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '389.44'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.48%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.95'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.18%  '
$ws.Range("E23").Value = '  +0.89%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.541'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.94%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.998'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.03%  '
$ws.Range("E27").Value = '  +0.42%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.38'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.52%  '
$ws.Range("E29").Value = '  -2.92%  '
$ws.Range("E30").Value = '  +0.13%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.23'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.04%  '
$ws.Range("E32").Value = '  -0.17%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '23.55'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.70%  '
$ws.Range("E35").Value = '  +0.39%  '
$ws.Range("E37").Value = '  -0.85%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '165.00'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.46%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.870'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.21%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.79'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +9.93%  '
$ws.Range("E41").Value = '  -1.21%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.81'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.83%  '
$ws.Range("E43").Value = '  +0.08%  '
$ws.Range("D44").Value = '2.842.67'
$ws.Range("E44").Value = '  +0.21%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '27.05'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.24%  '
$ws.Range("E46").Value = '  +0.31%  '
$ws.Range("E47").Value = '  -2.08%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '41.74'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.31%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0301'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.66%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '335.91'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.24%  '
$ws.Range("E51").Value = '  -1.47%  '
